$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E6").Value = "['Normal']"

$ws.Range("D12").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E12").Value = "['Normal', 'HardwareFault']"

$ws.Range("D15").Value = "[0, 0, 0, 0, 0, 0, 0]"
$ws.Range("E15").Value = "[]"

$ws.Range("D24").Value = "[0, 0, 0, 0, 0, 0, 0]"
$ws.Range("E24").Value = "[]"

$ws.Range("D28").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E28").Value = "['SoftwareFault']"

$ws.Range("D38").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E38").Value = "['Normal', 'HardwareFault']"

$ws.Range("D53").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E53").Value = "['Normal']"

$ws.Range("D56").Value = "[0, 0, 0, 0, 0, 0, 0]"
$ws.Range("E56").Value = "[]"

$ws.Range("D58").Value = "[0, 0, 0, 1, 0, 0, 0]"
$ws.Range("E58").Value = "['ParamViolation']"

$ws.Range("D61").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E61").Value = "['Normal', 'SoftwareFault']"

$ws.Range("D68").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E68").Value = "['Normal', 'ParamViolation']"

$ws.Range("D71").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E71").Value = "['Normal', 'ParamViolation']"

$ws.Range("D84").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E84").Value = "['Normal']"

$ws.Range("D88").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E88").Value = "['Normal']"

$ws.Range("D97").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E97").Value = "['Normal']"

$ws.Range("D113").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E113").Value = "['Normal', 'HardwareFault']"
